# Rename the second sheet "SAMPLE_TEST" -> "Exposure conditions"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SAMPLE_TEST")
$ws.Name = "Exposure conditions"

# Add new rows 6-11 with data (columns A-H stay empty, I/J/K/L populated)
$controlRows = 6..9
$i = 1
foreach ($r in $controlRows) {
    $ws.Cells.Item($r, 9).Value = $i
    $ws.Cells.Item($r, 10).Value = "CONTROL (SEE VEHICLE)"
    $ws.Cells.Item($r, 11).Value = 0
    $ws.Cells.Item($r, 12).Value = "TP1"
    $i++
}

$blankRows = 10..11
foreach ($r in $blankRows) {
    $ws.Cells.Item($r, 9).Value = 0
    $ws.Cells.Item($r, 10).Value = "EXTRACTION BLANK"
    $ws.Cells.Item($r, 11).NumberFormat = "@"
    $ws.Cells.Item($r, 11).Value = "0"
    $ws.Cells.Item($r, 12).Value = "TP0"
}
